$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Igfbp4", "Fzd8", "ECs", 3, 1, 6.112374333333332, 18.337123, 0.02125075796073207, 0.02125075796073207, 3, 1, 3.554362333333333, 10.663087, 0.2019894022634335, 0.2019894022634335, 21.72559309763344, 195.530337878701, 0.004292427898133172, 0.004292427898133173),
    @("ECs", "Igfbp4", "Fzd8", "FAPs", 3, 1, 6.112374333333332, 18.337123, 0.02125075796073207, 0.02125075796073207, 3, 1, 11.625787, 34.877361, 0.6606770910634029, 0.6606770910634029, 71.06116206360032, 639.5504585724029, 0.01403988895238892, 0.01403988895238892),
    @("ECs", "Igfbp4", "Fzd8", "M2", 3, 1, 6.112374333333332, 18.337123, 0.02125075796073207, 0.02125075796073207, 1, 0.3333333333333333, 0.02609533333333333, 0.078286, 0.001482961017348462, 0.001482961017348462, 0.1595044456864444, 1.435540011178, 0.00003151404564487315, 0.00003151404564487315),
    @("ECs", "Igfbp4", "Fzd8", "sCs", 3, 1, 6.112374333333332, 18.337123, 0.02125075796073207, 0.02125075796073207, 3, 1, 2.390531666666666, 7.171595, 0.1358505456558151, 0.1358505456558151, 14.61182440235389, 131.506419621185, 0.002886927064565108, 0.002886927064565109),
    @("FAPs", "Igfbp4", "Fzd8", "ECs", 3, 1, 263.0904693333333, 789.271408, 0.9146808720612395, 0.9146808720612395, 3, 1, 3.554362333333333, 10.663087, 0.2019894022634335, 0.2019894022634335, 935.1188544573884, 8416.069690116496, 0.1847558426094458, 0.1847558426094459),
    @("FAPs", "Igfbp4", "Fzd8", "FAPs", 3, 1, 263.0904693333333, 789.271408, 0.9146808720612395, 0.9146808720612395, 3, 1, 11.625787, 34.877361, 0.6606770910634029, 0.6606770910634029, 3058.633758199365, 27527.70382379429, 0.6043086978047563, 0.6043086978047563),
    @("FAPs", "Igfbp4", "Fzd8", "M2", 3, 1, 263.0904693333333, 789.271408, 0.9146808720612395, 0.9146808720612395, 1, 0.3333333333333333, 0.02609533333333333, 0.078286, 0.001482961017348462, 0.001482961017348462, 6.865433494076443, 61.78890144668799, 0.001356436076581114, 0.001356436076581114),
    @("FAPs", "Igfbp4", "Fzd8", "sCs", 3, 1, 263.0904693333333, 789.271408, 0.9146808720612395, 0.9146808720612395, 3, 1, 2.390531666666666, 7.171595, 0.1358505456558151, 0.1358505456558151, 628.9260981395288, 5660.334883255759, 0.1242598955704562, 0.1242598955704562),
    @("M2", "Igfbp4", "Fzd8", "ECs", 3, 1, 8.317515, 24.952545, 0.02891732221566466, 0.02891732221566466, 3, 1, 3.554362333333333, 10.663087, 0.2019894022634335, 0.2019894022634335, 29.563462022935, 266.071158206415, 0.00584099262940121, 0.005840992629401211),
    @("M2", "Igfbp4", "Fzd8", "FAPs", 3, 1, 8.317515, 24.952545, 0.02891732221566466, 0.02891732221566466, 3, 1, 11.625787, 34.877361, 0.6606770910634029, 0.6606770910634029, 96.697657759305, 870.278919833745, 0.01910501232278844, 0.01910501232278844),
    @("M2", "Igfbp4", "Fzd8", "M2", 3, 1, 8.317515, 24.952545, 0.02891732221566466, 0.02891732221566466, 1, 0.3333333333333333, 0.02609533333333333, 0.078286, 0.001482961017348462, 0.001482961017348462, 0.21704832643, 1.95343493787, 0.00004288326157193532, 0.00004288326157193532),
    @("M2", "Igfbp4", "Fzd8", "sCs", 3, 1, 8.317515, 24.952545, 0.02891732221566466, 0.02891732221566466, 3, 1, 2.390531666666666, 7.171595, 0.1358505456558151, 0.1358505456558151, 19.883282995475, 178.949546959275, 0.003928434001903067, 0.003928434001903068),
    @("sCs", "Igfbp4", "Fzd8", "ECs", 3, 1, 10.110527, 30.331581, 0.03515104776236379, 0.03515104776236379, 3, 1, 3.554362333333333, 10.663087, 0.2019894022634335, 0.2019894022634335, 35.93647633894966, 323.428287050547, 0.007100139126453264, 0.007100139126453265),
    @("sCs", "Igfbp4", "Fzd8", "FAPs", 3, 1, 10.110527, 30.331581, 0.03515104776236379, 0.03515104776236379, 3, 1, 11.625787, 34.877361, 0.6606770910634029, 0.6606770910634029, 117.542833359749, 1057.885500237741, 0.02322349198346925, 0.02322349198346925),
    @("sCs", "Igfbp4", "Fzd8", "M2", 3, 1, 10.110527, 30.331581, 0.03515104776236379, 0.03515104776236379, 1, 0.3333333333333333, 0.02609533333333333, 0.078286, 0.001482961017348462, 0.001482961017348462, 0.2638375722406666, 2.374538150166, 0.00005212763355053937, 0.00005212763355053937),
    @("sCs", "Igfbp4", "Fzd8", "sCs", 3, 1, 10.110527, 30.331581, 0.03515104776236379, 0.03515104776236379, 3, 1, 2.390531666666666, 7.171595, 0.1358505456558151, 0.1358505456558151, 24.16953496018833, 217.525814641695, 0.004775289018890739, 0.004775289018890741)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}
